# Apply the cryptos list update: refresh prices / volume(1h) figures,
# and correct the WEMIXTOKEN/Frax row ordering (rows 43-44).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") "28.431.13"
Set-TextValue $ws.Range("E2") "  +0.12%  "
Set-TextValue $ws.Range("D3") "1.831.54"
Set-TextValue $ws.Range("E3") "  +2.08%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.09%  "
Set-TextValue $ws.Range("D5") "317.88"
Set-TextValue $ws.Range("E5") "  +0.37%  "
Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  -0.04%  "
Set-TextValue $ws.Range("D7") "0.5305"
Set-TextValue $ws.Range("E7") "  -0.84%  "
Set-TextValue $ws.Range("D8") "0.4045"
Set-TextValue $ws.Range("E8") "  +7.47%  "
Set-TextValue $ws.Range("D9") "0.07565"
Set-TextValue $ws.Range("E9") "  +1.04%  "
Set-TextValue $ws.Range("D10") "41.99"
Set-TextValue $ws.Range("E10") "  +0.70%  "
Set-TextValue $ws.Range("E11") "  +0.79%  "
Set-TextValue $ws.Range("D12") "6.327"
Set-TextValue $ws.Range("E12") "  +3.23%  "
Set-TextValue $ws.Range("D13") "1.001"
Set-TextValue $ws.Range("E13") "  -0.21%  "
Set-TextValue $ws.Range("D14") "7.598"
Set-TextValue $ws.Range("E14") "  +4.32%  "
Set-TextValue $ws.Range("D15") "20.88"
Set-TextValue $ws.Range("E15") "  +1.46%  "
Set-TextValue $ws.Range("D16") "1.831.02"
Set-TextValue $ws.Range("E16") "  +1.41%  "
Set-TextValue $ws.Range("D17") "89.62"
Set-TextValue $ws.Range("E17") "  +0.17%  "
Set-TextValue $ws.Range("D18") "0.00001073"
Set-TextValue $ws.Range("E18") "  +1.18%  "
Set-TextValue $ws.Range("D19") "0.06616"
Set-TextValue $ws.Range("E19") "  +1.89%  "
Set-TextValue $ws.Range("E20") "  +1.16%  "
Set-TextValue $ws.Range("E21") "  -0.05%  "
Set-TextValue $ws.Range("D22") "6.067"
Set-TextValue $ws.Range("E22") "  +1.74%  "
Set-TextValue $ws.Range("D23") "28.463.44"
Set-TextValue $ws.Range("D24") "11.33"
Set-TextValue $ws.Range("E24") "  +2.28%  "
Set-TextValue $ws.Range("D25") "2.118"
Set-TextValue $ws.Range("E25") "  +1.72%  "
Set-TextValue $ws.Range("D26") "2.474"
Set-TextValue $ws.Range("E26") "  +7.78%  "
Set-TextValue $ws.Range("D27") "157.07"
Set-TextValue $ws.Range("E27") "  -1.49%  "
Set-TextValue $ws.Range("D28") "20.57"
Set-TextValue $ws.Range("E28") "  +1.10%  "
Set-TextValue $ws.Range("D29") "2.052.86"
Set-TextValue $ws.Range("E29") "  +2.54%  "
Set-TextValue $ws.Range("D30") "123.74"
Set-TextValue $ws.Range("E30") "  +1.11%  "
Set-TextValue $ws.Range("D31") "1.128"
Set-TextValue $ws.Range("E31") "  +2.66%  "
Set-TextValue $ws.Range("E32") "  +4.45%  "
Set-TextValue $ws.Range("D33") "5.696"
Set-TextValue $ws.Range("E33") "  +1.91%  "
Set-TextValue $ws.Range("D34") "3.657"
Set-TextValue $ws.Range("E34") "  -0.11%  "
Set-TextValue $ws.Range("D35") "0.07202"
Set-TextValue $ws.Range("E35") "  +9.91%  "
Set-TextValue $ws.Range("D36") "0.2273"
Set-TextValue $ws.Range("D37") "5.266"
Set-TextValue $ws.Range("E37") "  +5.44%  "
Set-TextValue $ws.Range("D38") "0.02350"
Set-TextValue $ws.Range("E38") "  +2.56%  "
Set-TextValue $ws.Range("D39") "8.816"
Set-TextValue $ws.Range("E39") "  +3.52%  "
Set-TextValue $ws.Range("E40") "  +2.56%  "
Set-TextValue $ws.Range("D41") "0.6269"
Set-TextValue $ws.Range("E41") "  +1.65%  "
Set-TextValue $ws.Range("D42") "1.193"
Set-TextValue $ws.Range("E42") "  +0.26%  "
Set-TextValue $ws.Range("B43") "WEMIXTOKEN"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D43") "1.413"
Set-TextValue $ws.Range("E43") "  -2.64%  "
Set-TextValue $ws.Range("B44") "Frax"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D44") "1.001"
Set-TextValue $ws.Range("E44") "  -0.04%  "
Set-TextValue $ws.Range("D45") "13.47"
Set-TextValue $ws.Range("E45") "  +1.60%  "
Set-TextValue $ws.Range("E46") "  +0.92%  "
Set-TextValue $ws.Range("D47") "0.5854"
Set-TextValue $ws.Range("E47") "  +1.17%  "
Set-TextValue $ws.Range("D48") "126.20"
Set-TextValue $ws.Range("E48") "  -0.23%  "
Set-TextValue $ws.Range("D49") "1.993"
Set-TextValue $ws.Range("E49") "  +2.94%  "
Set-TextValue $ws.Range("D50") "1.193"
Set-TextValue $ws.Range("E50") "  +0.14%  "
Set-TextValue $ws.Range("D51") "0.06906"
Set-TextValue $ws.Range("E51") "  +0.58%  "
